# Adds new functional "Pi" (Raspberry Pi) requirement statements to the
# "Functional Statements" sheet, and leaves the workbook focused on that
# sheet / the newly entered rows (matching the author's final view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional Statements")

# New requirement statements, typed in the order that assigns shared-string
# table indices 110-116 (matches the authoring order), before the 5th
# statement ("... from membership provider") is moved up to sit right after
# the other "overwrite ... with data from ..." statements.
$newStatements = @(
    "PI will be able to load configuration from SD card",
    "PI will be able to load initial membership data from SD card",
    "PI will be able to overwrite membership data on SD card with data from phone",
    "PI will be able to overwrite membership data on SD card with data from network",
    "PI will not allow access to SD card without special tools",
    "PI will not allow access to SD card from locked side (that is side with RFID scanner)",
    "PI will be able to overwrite membership data on SD card with data from membership provider"
)

$firstNewRow = 68
$r = $firstNewRow
foreach ($txt in $newStatements) {
    $ws.Cells.Item($r, 1).Value = "Pi"
    $ws.Cells.Item($r, 2).Value = $txt
    $r = $r + 1
}

# Re-order the last three new rows so the "membership provider" statement
# (typed last) ends up grouped with the other "overwrite ... with data from"
# statements, without re-typing the text (keeps shared-string allocation
# order intact).
$lastRow = $r - 1
$tmp = $ws.Range("B" + $lastRow).Value()
$ws.Range("B" + $lastRow).Value = $ws.Range("B" + ($lastRow - 1)).Value()
$ws.Range("B" + ($lastRow - 1)).Value = $ws.Range("B" + ($lastRow - 2)).Value()
$ws.Range("B" + ($lastRow - 2)).Value = $tmp

# Bring the Functional Statements sheet to the front and leave the
# selection on the last entered row, matching the saved view state.
$ws.Activate()
$ws.Range("B76").Select() | Out-Null
